$d = $word.ActiveDocument

# 1) Remove the stray "_GoBack" bookmark sitting before the "Operating
#    Conditions" heading.
$d.Bookmarks("_GoBack").Delete()

# 2) Find the "1000ms" value that follows the "Pass/De-mature Time" label
#    (this is the value that is being removed per the commit).
$valRng = $d.Content
$valRng.Find.Execute("Pass/De-mature Time")
$valRng.Collapse(0)
$valRng.Find.Execute("1000ms")
$valStart = $valRng.Start
$valEnd = $valRng.End

# 3) Re-create the "_GoBack" bookmark collapsed at the spot right after the
#    trailing tab (i.e. where the "1000ms" text currently starts). This has
#    to happen before the paragraph-merging delete in step 4, otherwise the
#    bookmark ends up re-anchored in the wrong place.
$bmRng = $d.Range($valStart, $valStart)
$d.Bookmarks.Add("_GoBack", $bmRng)

# 4) Delete the "1000ms" run(s).
$d.Range($valStart, $valEnd).Delete()

# 5) Delete the whole "Failure Aging ... 30 Ignition cycles" paragraph,
#    merging it away into the "Pass/De-mature Time" paragraph above it.
$failRng = $d.Content
$failRng.Find.Execute("Failure Aging")
$failPara = $failRng.Paragraphs(1)
$failPara.Range.Delete()
